# Auto-generated edit script: updates market-price derived columns (H:N)
# in the Ixion_Profits workbook across all 8 job sheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 9092650
$ws.Cells.Item(40, 9).Value = 1754.025
$ws.Cells.Item(40, 10).Value = 33335038
$ws.Cells.Item(40, 11).Value = 1754.025
$ws.Cells.Item(40, 12).Value = 33335038
$ws.Cells.Item(40, 13).Value = -1579.025
$ws.Cells.Item(40, 14).Value = -33335388
$ws.Cells.Item(80, 8).Value = 5694.7896
$ws.Cells.Item(80, 9).Value = 369.23077
$ws.Cells.Item(80, 10).Value = 17233.5
$ws.Cells.Item(80, 11).Value = 1107.69231
$ws.Cells.Item(80, 12).Value = 51700.5
$ws.Cells.Item(80, 13).Value = -109.6923099999999
$ws.Cells.Item(80, 14).Value = -53696.5
$ws.Cells.Item(83, 8).Value = 5694.7896
$ws.Cells.Item(83, 9).Value = 369.23077
$ws.Cells.Item(83, 10).Value = 17233.5
$ws.Cells.Item(83, 11).Value = 3323.07693
$ws.Cells.Item(83, 12).Value = 155101.5
$ws.Cells.Item(83, 13).Value = 1668.92307
$ws.Cells.Item(83, 14).Value = -165085.5
$ws.Cells.Item(92, 8).Value = 38580630
$ws.Cells.Item(92, 9).Value = 1792438
$ws.Cells.Item(92, 11).Value = 1792438
$ws.Cells.Item(92, 13).Value = -1791190
$ws.Cells.Item(98, 8).Value = 338.22223
$ws.Cells.Item(98, 9).Value = 342.25
$ws.Cells.Item(98, 10).Value = 306
$ws.Cells.Item(98, 11).Value = 342.25
$ws.Cells.Item(98, 12).Value = 306
$ws.Cells.Item(98, 13).Value = 1155.75
$ws.Cells.Item(98, 14).Value = -3302
$ws.Cells.Item(112, 8).Value = 17317282
$ws.Cells.Item(112, 9).Value = 850
$ws.Cells.Item(112, 10).Value = 19048926
$ws.Cells.Item(112, 11).Value = 2550
$ws.Cells.Item(112, 12).Value = 57146778
$ws.Cells.Item(112, 13).Value = -1442
$ws.Cells.Item(112, 14).Value = -57148994
$ws.Cells.Item(122, 8).Value = 338.22223
$ws.Cells.Item(122, 9).Value = 342.25
$ws.Cells.Item(122, 10).Value = 306
$ws.Cells.Item(122, 11).Value = 1026.75
$ws.Cells.Item(122, 12).Value = 918
$ws.Cells.Item(122, 13).Value = 1423.25
$ws.Cells.Item(122, 14).Value = -5818
$ws.Cells.Item(129, 8).Value = 992.73334
$ws.Cells.Item(129, 9).Value = 635.2105
$ws.Cells.Item(129, 10).Value = 1158.4147
$ws.Cells.Item(129, 11).Value = 1905.6315
$ws.Cells.Item(129, 12).Value = 3475.2441
$ws.Cells.Item(129, 13).Value = 3094.3685
$ws.Cells.Item(129, 14).Value = -13475.2441
$ws.Cells.Item(131, 8).Value = 2104.45
$ws.Cells.Item(131, 10).Value = 2936
$ws.Cells.Item(131, 12).Value = 8808
$ws.Cells.Item(131, 14).Value = -18888
$ws.Cells.Item(132, 8).Value = 1691.174
$ws.Cells.Item(132, 9).Value = 1291.9531
$ws.Cells.Item(132, 11).Value = 3875.8593
$ws.Cells.Item(132, 13).Value = -1345.8593
$ws.Cells.Item(137, 8).Value = 1805.8235
$ws.Cells.Item(137, 9).Value = 2079.3
$ws.Cells.Item(137, 10).Value = 1415.1428
$ws.Cells.Item(137, 11).Value = 6237.900000000001
$ws.Cells.Item(137, 12).Value = 4245.428400000001
$ws.Cells.Item(137, 13).Value = -3687.900000000001
$ws.Cells.Item(137, 14).Value = -9345.428400000001
$ws.Cells.Item(138, 8).Value = 4426.8447
$ws.Cells.Item(138, 9).Value = 2032.7894
$ws.Cells.Item(138, 10).Value = 5593.1797
$ws.Cells.Item(138, 11).Value = 6098.3682
$ws.Cells.Item(138, 12).Value = 16779.5391
$ws.Cells.Item(138, 13).Value = -958.3681999999999
$ws.Cells.Item(138, 14).Value = -27059.5391

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 19037.74
$ws.Cells.Item(32, 9).Value = 16183.14
$ws.Cells.Item(32, 10).Value = 36573.145
$ws.Cells.Item(32, 11).Value = 16183.14
$ws.Cells.Item(32, 12).Value = 36573.145
$ws.Cells.Item(32, 13).Value = -15896.14
$ws.Cells.Item(32, 14).Value = -37147.145
$ws.Cells.Item(45, 8).Value = 176410.33
$ws.Cells.Item(45, 9).Value = 300632
$ws.Cells.Item(45, 10).Value = 2500
$ws.Cells.Item(45, 11).Value = 300632
$ws.Cells.Item(45, 12).Value = 2500
$ws.Cells.Item(45, 13).Value = -300255
$ws.Cells.Item(45, 14).Value = -3254
$ws.Cells.Item(97, 8).Value = 611.4
$ws.Cells.Item(97, 9).Value = 481.53845
$ws.Cells.Item(97, 10).Value = 1455.5
$ws.Cells.Item(97, 11).Value = 481.53845
$ws.Cells.Item(97, 12).Value = 1455.5
$ws.Cells.Item(97, 13).Value = 14.46154999999999
$ws.Cells.Item(97, 14).Value = -2447.5
$ws.Cells.Item(102, 8).Value = 2646755.8
$ws.Cells.Item(102, 9).Value = 2646755.8
$ws.Cells.Item(102, 11).Value = 2646755.8
$ws.Cells.Item(102, 13).Value = -2645133.8
$ws.Cells.Item(132, 8).Value = 3724.5405
$ws.Cells.Item(132, 9).Value = 2794.889
$ws.Cells.Item(132, 10).Value = 4605.263
$ws.Cells.Item(132, 11).Value = 8384.667000000001
$ws.Cells.Item(132, 12).Value = 13815.789
$ws.Cells.Item(132, 13).Value = -5854.667000000001
$ws.Cells.Item(132, 14).Value = -18875.789

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 337.48
$ws.Cells.Item(80, 9).Value = 78.333336
$ws.Cells.Item(80, 10).Value = 419.3158
$ws.Cells.Item(80, 11).Value = 78.333336
$ws.Cells.Item(80, 12).Value = 419.3158
$ws.Cells.Item(80, 13).Value = 919.666664
$ws.Cells.Item(80, 14).Value = -2415.3158
$ws.Cells.Item(83, 8).Value = 337.48
$ws.Cells.Item(83, 9).Value = 78.333336
$ws.Cells.Item(83, 10).Value = 419.3158
$ws.Cells.Item(83, 11).Value = 391.66668
$ws.Cells.Item(83, 12).Value = 2096.579
$ws.Cells.Item(83, 13).Value = 4600.33332
$ws.Cells.Item(83, 14).Value = -12080.579
$ws.Cells.Item(99, 8).Value = 31251528
$ws.Cells.Item(99, 9).Value = 43479652
$ws.Cells.Item(99, 11).Value = 43479652
$ws.Cells.Item(99, 13).Value = -43478154
$ws.Cells.Item(134, 8).Value = 60616.79
$ws.Cells.Item(134, 9).Value = 8428.056
$ws.Cells.Item(134, 10).Value = 1000014
$ws.Cells.Item(134, 11).Value = 25284.168
$ws.Cells.Item(134, 12).Value = 3000042
$ws.Cells.Item(134, 13).Value = -22749.168
$ws.Cells.Item(134, 14).Value = -3005112

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 286.83334
$ws.Cells.Item(22, 9).Value = 127.28571
$ws.Cells.Item(22, 11).Value = 127.28571
$ws.Cells.Item(22, 13).Value = 222.71429
$ws.Cells.Item(31, 8).Value = 4740.963
$ws.Cells.Item(31, 9).Value = 1510.8235
$ws.Cells.Item(31, 10).Value = 10232.2
$ws.Cells.Item(31, 11).Value = 1510.8235
$ws.Cells.Item(31, 12).Value = 10232.2
$ws.Cells.Item(31, 13).Value = -1215.8235
$ws.Cells.Item(31, 14).Value = -10822.2
$ws.Cells.Item(34, 8).Value = 4740.963
$ws.Cells.Item(34, 9).Value = 1510.8235
$ws.Cells.Item(34, 10).Value = 10232.2
$ws.Cells.Item(34, 11).Value = 1510.8235
$ws.Cells.Item(34, 12).Value = 10232.2
$ws.Cells.Item(34, 13).Value = -1308.8235
$ws.Cells.Item(34, 14).Value = -10636.2
$ws.Cells.Item(62, 8).Value = 5080
$ws.Cells.Item(62, 9).Value = 6032
$ws.Cells.Item(62, 11).Value = 6032
$ws.Cells.Item(62, 13).Value = -5408
$ws.Cells.Item(65, 8).Value = 5080
$ws.Cells.Item(65, 9).Value = 6032
$ws.Cells.Item(65, 11).Value = 30160
$ws.Cells.Item(65, 13).Value = -27040
$ws.Cells.Item(68, 8).Value = 32000
$ws.Cells.Item(68, 10).Value = 32000
$ws.Cells.Item(68, 12).Value = 32000
$ws.Cells.Item(68, 14).Value = -33498
$ws.Cells.Item(71, 8).Value = 32000
$ws.Cells.Item(71, 10).Value = 32000
$ws.Cells.Item(71, 12).Value = 96000
$ws.Cells.Item(71, 14).Value = -103488
$ws.Cells.Item(105, 8).Value = 2685.0625
$ws.Cells.Item(105, 9).Value = 2589.2856
$ws.Cells.Item(105, 11).Value = 2589.2856
$ws.Cells.Item(105, 13).Value = -842.2856000000002
$ws.Cells.Item(122, 8).Value = 4509.5293
$ws.Cells.Item(122, 9).Value = 3546.5454
$ws.Cells.Item(122, 11).Value = 10639.6362
$ws.Cells.Item(122, 13).Value = -8189.636200000001
$ws.Cells.Item(132, 8).Value = 2343.36
$ws.Cells.Item(132, 9).Value = 1756
$ws.Cells.Item(132, 10).Value = 3853.7144
$ws.Cells.Item(132, 11).Value = 5268
$ws.Cells.Item(132, 12).Value = 11561.1432
$ws.Cells.Item(132, 13).Value = -2738
$ws.Cells.Item(132, 14).Value = -16621.1432
$ws.Cells.Item(134, 8).Value = 288687.34
$ws.Cells.Item(134, 9).Value = 3607.577
$ws.Cells.Item(134, 10).Value = 1112251.1
$ws.Cells.Item(134, 11).Value = 10822.731
$ws.Cells.Item(134, 12).Value = 3336753.3
$ws.Cells.Item(134, 13).Value = -8287.731
$ws.Cells.Item(134, 14).Value = -3341823.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 213261.11
$ws.Cells.Item(113, 10).Value = 286205.97
$ws.Cells.Item(113, 12).Value = 858617.9099999999
$ws.Cells.Item(113, 14).Value = -862957.9099999999
$ws.Cells.Item(122, 8).Value = 6878.4443
$ws.Cells.Item(122, 10).Value = 15932.714
$ws.Cells.Item(122, 12).Value = 143394.426
$ws.Cells.Item(122, 14).Value = -148294.426
$ws.Cells.Item(131, 8).Value = 3449329.2
$ws.Cells.Item(131, 9).Value = 14286519
$ws.Cells.Item(131, 10).Value = 1132.6364
$ws.Cells.Item(131, 11).Value = 42859557
$ws.Cells.Item(131, 12).Value = 3397.9092
$ws.Cells.Item(131, 13).Value = -42854517
$ws.Cells.Item(131, 14).Value = -13477.9092
$ws.Cells.Item(132, 8).Value = 3105.9048
$ws.Cells.Item(132, 10).Value = 3242.353
$ws.Cells.Item(132, 12).Value = 29181.177
$ws.Cells.Item(132, 14).Value = -34241.177

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 10526.333
$ws.Cells.Item(80, 9).Value = 26602.5
$ws.Cells.Item(80, 10).Value = 2488.25
$ws.Cells.Item(80, 11).Value = 26602.5
$ws.Cells.Item(80, 12).Value = 2488.25
$ws.Cells.Item(80, 13).Value = -25604.5
$ws.Cells.Item(80, 14).Value = -4484.25
$ws.Cells.Item(83, 8).Value = 10526.333
$ws.Cells.Item(83, 9).Value = 26602.5
$ws.Cells.Item(83, 10).Value = 2488.25
$ws.Cells.Item(83, 11).Value = 133012.5
$ws.Cells.Item(83, 12).Value = 12441.25
$ws.Cells.Item(83, 13).Value = -128020.5
$ws.Cells.Item(83, 14).Value = -22425.25
$ws.Cells.Item(107, 8).Value = 765.7727
$ws.Cells.Item(107, 9).Value = 604.63635
$ws.Cells.Item(107, 11).Value = 604.63635
$ws.Cells.Item(107, 13).Value = 1315.36365
$ws.Cells.Item(113, 8).Value = 76931710
$ws.Cells.Item(113, 9).Value = 100010840
$ws.Cells.Item(113, 10).Value = 1300
$ws.Cells.Item(113, 11).Value = 100010840
$ws.Cells.Item(113, 12).Value = 1300
$ws.Cells.Item(113, 13).Value = -100008670
$ws.Cells.Item(113, 14).Value = -5640
$ws.Cells.Item(132, 8).Value = 2169.1052
$ws.Cells.Item(132, 9).Value = 1894.9412
$ws.Cells.Item(132, 10).Value = 4499.5
$ws.Cells.Item(132, 11).Value = 5684.8236
$ws.Cells.Item(132, 12).Value = 13498.5
$ws.Cells.Item(132, 13).Value = -3154.8236
$ws.Cells.Item(132, 14).Value = -18558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1724.35
$ws.Cells.Item(100, 9).Value = 1406.4375
$ws.Cells.Item(100, 11).Value = 1406.4375
$ws.Cells.Item(100, 13).Value = -865.4375
$ws.Cells.Item(136, 8).Value = 9411.75
$ws.Cells.Item(136, 9).Value = 6222.846
$ws.Cells.Item(136, 10).Value = 17702.9
$ws.Cells.Item(136, 11).Value = 18668.538
$ws.Cells.Item(136, 12).Value = 53108.7
$ws.Cells.Item(136, 13).Value = -16118.538
$ws.Cells.Item(136, 14).Value = -58208.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value = 32875
$ws.Cells.Item(123, 10).Value = 32875
$ws.Cells.Item(123, 12).Value = 32875
$ws.Cells.Item(123, 14).Value = -42675
$ws.Cells.Item(132, 8).Value = 1792.3793
$ws.Cells.Item(132, 9).Value = 933.7273
$ws.Cells.Item(132, 10).Value = 2317.111
$ws.Cells.Item(132, 11).Value = 2801.1819
$ws.Cells.Item(132, 12).Value = 6951.333
$ws.Cells.Item(132, 13).Value = -271.1819
$ws.Cells.Item(132, 14).Value = -12011.333

